# Apply crypto price/volume updates scraped on Sun Mar 17 13:52:18 UTC 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.034.83"
$ws.Range("E2").Value = "  -1.44%  "

# Row 3
$ws.Range("D3").Value = "3.565.56"
$ws.Range("E3").Value = "  -2.76%  "

# Row 4
$ws.Range("E4").Value = "  +0.09%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "191.59"
$ws.Range("E5").Value = "  -0.35%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "570.60"
$ws.Range("E6").Value = "  -4.45%  "

# Row 7
$ws.Range("D7").Value = "3.558.96"
$ws.Range("E7").Value = "  -2.80%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.614"
$ws.Range("E8").Value = "  -0.82%  "

# Row 9
$ws.Range("E9").Value = "  +0.27%  "

# Row 10
$ws.Range("E10").Value = "  -4.02%  "

# Row 11
$ws.Range("B11").Value = "Dogecoin"
$ws.Range("C11").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("E11").Value = "  -2.22%  "

# Row 12
$ws.Range("B12").Value = "Avalanche"
$ws.Range("C12").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "55.50"
$ws.Range("E12").Value = "  -3.33%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000270"
$ws.Range("E13").Value = "  -0.63%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.85"
$ws.Range("E14").Value = "  -3.36%  "

# Row 15
$ws.Range("D15").Value = "4.132.19"
$ws.Range("E15").Value = "  -2.83%  "

# Row 16
$ws.Range("D16").Value = "3.564.45"
$ws.Range("E16").Value = "  -2.86%  "

# Row 17
$ws.Range("E17").Value = "  -1.17%  "

# Row 18
$ws.Range("D18").Value = "66.868.16"
$ws.Range("E18").Value = "  -1.38%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.18"
$ws.Range("E19").Value = "  -2.14%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.17"
$ws.Range("E20").Value = "  -3.72%  "

# Row 21
$ws.Range("E21").Value = "  -4.95%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "401.38"
$ws.Range("E22").Value = "  +0.40%  "

# Row 23
$ws.Range("E23").Value = "  -6.13%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.90"
$ws.Range("E24").Value = "  +5.88%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.57"
$ws.Range("E25").Value = "  -2.22%  "

# Row 26
$ws.Range("E26").Value = "  -1.66%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.47"
$ws.Range("E27").Value = "  +0.11%  "

# Row 28
$ws.Range("E28").Value = "  +0.97%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "3.69"
$ws.Range("E29").Value = "  +0.95%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.80"
$ws.Range("E30").Value = "  +5.63%  "

# Row 31
$ws.Range("E31").Value = "  -3.66%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "31.14"
$ws.Range("E32").Value = "  -2.11%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "643.28"
$ws.Range("E33").Value = "  +5.44%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "12.10"
$ws.Range("E34").Value = "  -1.55%  "

# Row 35
$ws.Range("E35").Value = "  -2.43%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "63.82"
$ws.Range("E36").Value = "  -5.89%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "42.24"
$ws.Range("E37").Value = "  -6.24%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.407"
$ws.Range("E38").Value = "  +3.19%  "

# Row 39
$ws.Range("E39").Value = "  +0.27%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0764"
$ws.Range("E40").Value = "  +0.60%  "

# Row 41
$ws.Range("B41").Value = "ThetaToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.11"
$ws.Range("E41").Value = "  +7.23%  "

# Row 42
$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "3.175.18"
$ws.Range("E42").Value = "  +13.02%  "

# Row 43
$ws.Range("B43").Value = "Kaspa"
$ws.Range("C43").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.133"
$ws.Range("E43").Value = "  -0.71%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.999"
$ws.Range("E44").Value = "  -0.03%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.70"
$ws.Range("E45").Value = "  +6.45%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0413"
$ws.Range("E46").Value = "  -2.80%  "

# Row 47
$ws.Range("E47").Value = "  -1.82%  "

# Row 48
$ws.Range("B48").Value = "Monero"
$ws.Range("C48").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "142.84"
$ws.Range("E48").Value = "  -0.73%  "

# Row 49
$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.130"
$ws.Range("E49").Value = "  -3.92%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.53"
$ws.Range("E50").Value = "  -4.35%  "

# Row 51
$ws.Range("E51").Value = "  -4.10%  "
